# "made Team as String" - populate the bracket template placeholders
# (tournament name / birth-year range / weight category / participant
# names+team) as literal FreeMarker-style text, upgrade the "weight
# category" cell to a bold Arial-10 left-aligned style, add the
# "jx:area(lastCell="P46")" threaded comment on A1, and move the active
# selection to T9 (matching the saved workbook state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header placeholders -----------------------------------------------
$ws.Range("C1").Value = '${tournamentName}'
$ws.Range("O3").Value = '${birthYearRange}'

# "Team" -> String: weight category becomes a bold/left-aligned Arial 10
# cell (new font/style), rather than the previous unstyled bold-8 cell.
$ws.Range("O5").Value = '${weightCategory}'
$ws.Range("O5").Font.Bold = $true
$ws.Range("O5").Font.Name = "Arial"
$ws.Range("O5").Font.Size = 10
$ws.Range("O5").HorizontalAlignment = -4131

# --- Bracket participant name/team placeholders -------------------------
$ws.Range("B12").Value = '${graph.get("THREE").get(7).participant.lastName} ${graph.get("THREE").get(7).participant.firstName} ${graph.get("THREE").get(7).participant.team? "(" + graph.get("THREE").get(7).participant.team + ")" : null}'
$ws.Range("B14").Value = '${graph.get("THREE").get(6).participant.lastName} ${graph.get("THREE").get(6).participant.firstName} ${graph.get("THREE").get(6).participant.team? "(" + graph.get("THREE").get(6).participant.team + ")" : null}'
$ws.Range("B16").Value = '${graph.get("THREE").get(5).participant.lastName} ${graph.get("THREE").get(5).participant.firstName} ${graph.get("THREE").get(5).participant.team? "(" + graph.get("THREE").get(5).participant.team + ")" : null}'
$ws.Range("B18").Value = '${graph.get("THREE").get(4).participant.lastName} ${graph.get("THREE").get(4).participant.firstName} ${graph.get("THREE").get(4).participant.team? "(" + graph.get("THREE").get(4).participant.team + ")" : null}'
$ws.Range("B22").Value = '${graph.get("THREE").get(3).participant.lastName} ${graph.get("THREE").get(3).participant.firstName} ${graph.get("THREE").get(3).participant.team? "(" + graph.get("THREE").get(3).participant.team + ")" : null}'
$ws.Range("B24").Value = '${graph.get("THREE").get(2).participant.lastName} ${graph.get("THREE").get(2).participant.firstName} ${graph.get("THREE").get(2).participant.team? "(" + graph.get("THREE").get(2).participant.team + ")" : null}'
$ws.Range("B26").Value = '${graph.get("THREE").get(1).participant.lastName} ${graph.get("THREE").get(1).participant.firstName} ${graph.get("THREE").get(1).participant.team? "(" + graph.get("THREE").get(1).participant.team + ")" : null}'
$ws.Range("B28").Value = '${graph.get("THREE").get(0).participant.lastName} ${graph.get("THREE").get(0).participant.firstName} ${graph.get("THREE").get(0).participant.team? "(" + graph.get("THREE").get(0).participant.team + ")" : null}'

# NB: shared-string table order follows insertion order, and the saved
# workbook built the "TWO" bracket strings bottom-up (E27 first ... E13
# last), so replicate that order here to land on matching sharedStrings
# indices.
$ws.Range("E27").Value = '${graph.get("TWO").get(0).participant.lastName} ${graph.get("TWO").get(0).participant.firstName} ${graph.get("TWO").get(0).participant.team? "(" + graph.get("TWO").get(0).participant.team + ")" : null}'
$ws.Range("E23").Value = '${graph.get("TWO").get(1).participant.lastName} ${graph.get("TWO").get(1).participant.firstName} ${graph.get("TWO").get(1).participant.team? "(" + graph.get("TWO").get(1).participant.team + ")" : null}'
$ws.Range("E17").Value = '${graph.get("TWO").get(2).participant.lastName} ${graph.get("TWO").get(2).participant.firstName} ${graph.get("TWO").get(2).participant.team? "(" + graph.get("TWO").get(2).participant.team + ")" : null}'
$ws.Range("E13").Value = '${graph.get("TWO").get(3).participant.lastName} ${graph.get("TWO").get(3).participant.firstName} ${graph.get("TWO").get(3).participant.team? "(" + graph.get("TWO").get(3).participant.team + ")" : null}'

# --- Threaded comment defining the print/export area --------------------
$ws.Range("A1").AddCommentThreaded('jx:area(lastCell="P46")') | Out-Null

# --- Selection state ------------------------------------------------------
$ws.Range("T9").Select() | Out-Null
